$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("G2").Value = 0.001364
$ws.Range("H2").Value = 0.004092
$ws.Range("M2").Value = 0.04257266666666667
$ws.Range("N2").Value = 0.127718
$ws.Range("O2").Value = 0.007034104319295177
$ws.Range("P2").Value = 0.007034104319295176
$ws.Range("Q2").Value = [double]"5.806911733333333E-05"
$ws.Range("R2").Value = 0.000522622056
$ws.Range("S2").Value = 0.007034104319295177
$ws.Range("T2").Value = 0.007034104319295176

# Update row 3 values
$ws.Range("G3").Value = 0.001364
$ws.Range("H3").Value = 0.004092
$ws.Range("M3").Value = 0.04484833333333333
$ws.Range("N3").Value = 0.134545
$ws.Range("O3").Value = 0.007410103240260335
$ws.Range("P3").Value = 0.007410103240260335
$ws.Range("Q3").Value = [double]"6.117312666666666E-05"
$ws.Range("R3").Value = 0.00055055814
$ws.Range("S3").Value = 0.007410103240260335
$ws.Range("T3").Value = 0.007410103240260335

# Update row 4 values
$ws.Range("G4").Value = 0.001364
$ws.Range("H4").Value = 0.004092
$ws.Range("M4").Value = 5.964901333333334
$ws.Range("N4").Value = 17.894704
$ws.Range("O4").Value = 0.9855557924404444
$ws.Range("P4").Value = 0.9855557924404444
$ws.Range("Q4").Value = 0.008136125418666667
$ws.Range("R4").Value = 0.07322512876800001
$ws.Range("S4").Value = 0.9855557924404444
$ws.Range("T4").Value = 0.9855557924404444

# Delete row 5 entirely (Resolving-Mac row), shifting cells up
$ws.Rows.Item(5).Delete()
